$wb = $excel.ActiveWorkbook

# --- Sheet "Obras" (sheet1): insert ID column at front, add Referentes column at end ---
$ws1 = $wb.Worksheets.Item("Obras")

# Insert a new column before column A; this shifts the existing A:E data to B:F.
$ws1.Columns.Item(1).Insert()

# New header cells: ID (col A) and Referentes (col G)
$ws1.Range("A1").Value = "ID"
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("G1").Value = "Referentes"
$ws1.Range("G1").Font.Bold = $true

# ID values for each data row
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3

# Referentes values (relation counts/ids) for each data row
$ws1.Range("G2").Value = 1
$ws1.Range("G3").Value = 2.3
$ws1.Range("G4").Value = 4

# Archivo column (now F) was manually widened for the new layout
$ws1.Columns.Item(6).ColumnWidth = 17.6

$ws1.Range("G4").Select() | Out-Null

# --- Sheet "Referentes" (sheet2): insert ID column at front ---
$ws2 = $wb.Worksheets.Item("Referentes")

# Insert a new column before column A; this shifts the existing A:D data to B:E.
$ws2.Columns.Item(1).Insert()

$ws2.Range("A1").Value = "ID"

# ID values for each data row
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
